$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename "Inference Time" header (C11) to a more descriptive label.
$ws.Range("C11").Value = "Inference Time s/traj (colab pro+ CPU)"

# Update EKF inference-time formula (C14) with the newly measured value.
$ws.Range("C14").Formula = "=2.12102031707763/10"

# Fill in newly measured EKF std inference time (B15).
$ws.Range("B15").Value = 0.63129999999999997

# Update RTS inference-time formula (C16) with the newly measured value.
$ws.Range("C16").Formula = "=4.10807561874389/10"

# Fill in newly measured RTS std inference time (B17).
$ws.Range("B17").Value = 0.89549999999999996

# Fill in newly measured RTSNet [dB] result (B18) and its inference time (C18).
$ws.Range("B18").Value = -23.454699999999999
$ws.Range("C18").Formula = "=4.35626101493835/10"

# Fill in newly measured RTSNet std inference time (B19).
$ws.Range("B19").Value = 0.69079999999999997

# Update the sheet's selection, matching where the author left off editing.
$ws.Range("D18").Select()
